$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2026-01-29 Thursday" "2026-01-30 Friday"

Replace-Text "632÷4=158, 0" "204÷5=40, 4"
Replace-Text "324÷9=36, 0" "195÷5=39, 0"
Replace-Text "175÷7=25, 0" "238÷9=26, 4"
Replace-Text "344÷8=43, 0" "319÷8=39, 7"
Replace-Text "175÷6=29, 1" "264÷6=44, 0"

Replace-Text "853÷8=106, 5" "911÷8=113, 7"
Replace-Text "233÷5=46, 3" "867÷8=108, 3"
Replace-Text "451÷9=50, 1" "201÷5=40, 1"
Replace-Text "422÷7=60, 2" "376÷3=125, 1"
Replace-Text "668÷4=167, 0" "429÷6=71, 3"

Replace-Text "865÷3=288, 1" "463÷6=77, 1"
Replace-Text "697÷7=99, 4" "745÷5=149, 0"
Replace-Text "494÷5=98, 4" "344÷7=49, 1"
Replace-Text "231÷9=25, 6" "314÷9=34, 8"
Replace-Text "624÷7=89, 1" "620÷2=310, 0"

Replace-Text "220÷3=73, 1" "976÷4=244, 0"
Replace-Text "775÷3=258, 1" "581÷6=96, 5"
Replace-Text "120÷8=15, 0" "454÷3=151, 1"
Replace-Text "100÷3=33, 1" "382÷3=127, 1"
Replace-Text "974÷8=121, 6" "243÷3=81, 0"

Replace-Text "152÷2=76, 0" "939÷5=187, 4"
Replace-Text "566÷5=113, 1" "376÷3=125, 1"
Replace-Text "826÷2=413, 0" "558÷9=62, 0"
Replace-Text "248÷2=124, 0" "993÷3=331, 0"
Replace-Text "620÷8=77, 4" "632÷7=90, 2"

Write-Output "Done"
